# Weekly data update: insert a new week's price group (Especial/Primera/Segunda)
# for Mango at the top of the reverse-chronological data table (before the
# row currently holding 2021-05-05 / serial 44321), pushing all subsequent
# rows down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 417 (formats get inherited from the row above,
# which already carries the date style for column D).
$ws.Range("A417:A419").EntireRow.Insert()

# Row 417 - Especial
$ws.Cells.Item(417, 1).Value = 8
$ws.Cells.Item(417, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(417, 3).Value = "Coquimbo"
$ws.Cells.Item(417, 4).Value = 44543
$ws.Cells.Item(417, 5).Value = 4
$ws.Cells.Item(417, 6).Value = "Fruta"
$ws.Cells.Item(417, 7).Value = 100108
$ws.Cells.Item(417, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(417, 9).Value = 100108002
$ws.Cells.Item(417, 10).Value = "Mango"
$ws.Cells.Item(417, 11).Value = "Sin especificar"
$ws.Cells.Item(417, 12).Value = "Especial"
$ws.Cells.Item(417, 13).Value = 512
$ws.Cells.Item(417, 14).Value = 5500
$ws.Cells.Item(417, 15).Value = 6000
$ws.Cells.Item(417, 16).Value = 5750
$ws.Cells.Item(417, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(417, 18).Value = "Perú"
$ws.Cells.Item(417, 19).Value = 1438
$ws.Cells.Item(417, 20).Value = 4

# Row 418 - Primera
$ws.Cells.Item(418, 1).Value = 8
$ws.Cells.Item(418, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(418, 3).Value = "Coquimbo"
$ws.Cells.Item(418, 4).Value = 44543
$ws.Cells.Item(418, 5).Value = 4
$ws.Cells.Item(418, 6).Value = "Fruta"
$ws.Cells.Item(418, 7).Value = 100108
$ws.Cells.Item(418, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(418, 9).Value = 100108002
$ws.Cells.Item(418, 10).Value = "Mango"
$ws.Cells.Item(418, 11).Value = "Sin especificar"
$ws.Cells.Item(418, 12).Value = "Primera"
$ws.Cells.Item(418, 13).Value = 512
$ws.Cells.Item(418, 14).Value = 5500
$ws.Cells.Item(418, 15).Value = 6000
$ws.Cells.Item(418, 16).Value = 5750
$ws.Cells.Item(418, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(418, 18).Value = "Perú"
$ws.Cells.Item(418, 19).Value = 1438
$ws.Cells.Item(418, 20).Value = 4

# Row 419 - Segunda
$ws.Cells.Item(419, 1).Value = 8
$ws.Cells.Item(419, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(419, 3).Value = "Coquimbo"
$ws.Cells.Item(419, 4).Value = 44543
$ws.Cells.Item(419, 5).Value = 4
$ws.Cells.Item(419, 6).Value = "Fruta"
$ws.Cells.Item(419, 7).Value = 100108
$ws.Cells.Item(419, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(419, 9).Value = 100108002
$ws.Cells.Item(419, 10).Value = "Mango"
$ws.Cells.Item(419, 11).Value = "Sin especificar"
$ws.Cells.Item(419, 12).Value = "Segunda"
$ws.Cells.Item(419, 13).Value = 512
$ws.Cells.Item(419, 14).Value = 5500
$ws.Cells.Item(419, 15).Value = 6000
$ws.Cells.Item(419, 16).Value = 5750
$ws.Cells.Item(419, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(419, 18).Value = "Perú"
$ws.Cells.Item(419, 19).Value = 1438
$ws.Cells.Item(419, 20).Value = 4
